$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.999.61"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.785.75"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.56"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.09%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5362"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.67%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3763"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07449"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.05%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.70"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.16%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.093"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("E12").Value = "  +0.08%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.55"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.92%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.096"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.73%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.209"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").Value = "1.778.50"
$ws.Range("E16").Value = "  -1.39%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "88.71"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.76%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06449"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.29"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.63%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.895"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "28.020.12"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  -2.36%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.087"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.14%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "155.31"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.18%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.25"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "1.982.07"
$ws.Range("E28").Value = "  -1.54%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.281"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.08%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "120.08"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.17%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.109"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("E32").Value = "  +3.14%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.635"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.544"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.86%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.2257"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.97%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.06437"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.86%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02284"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.92%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.013"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.85%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.463"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6153"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.444"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.26%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "11.08"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.95%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.171"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.22"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.667"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5756"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "126.81"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.187"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.39%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.924"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.44%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06796"
$c.Style = "Normal"
